# Atualizei dados bibi e add
# Update faturamento_anual data for row 9 (Ano = 2025)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B9").Value = 3873874.99
$ws.Range("C9").Value = 613025.24
$ws.Range("D9").Value = 4486900.23
$ws.Range("E9").Value = 13.66255563030426
$ws.Range("F9").Value = 86.33744436969573
$ws.Range("G9").Value = -40.75416723876488
$ws.Range("H9").Value = -30.04317230407292
$ws.Range("I9").Value = 39365
$ws.Range("J9").Value = 1683
$ws.Range("K9").Value = 41048
$ws.Range("L9").Value = 28383
$ws.Range("M9").Value = 158.0840725081915
$ws.Range("N9").Value = 7.927216114955193
